$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "section" column (A) for rows that belong to section 1
$ws.Range("A4").Value = 1
$ws.Range("A5").Value = 1
$ws.Range("A6").Value = 1

# Fill in the "section" column (A) for rows that belong to section 2
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 2
$ws.Range("A10").Value = 2

# Update the last selected cell on the sheet
$ws.Range("C13").Select()
